# Update the individual VIN values used by the SS tests (pt2).
# Shared VIN string "1GPGP1111&1" -> "AAAKN3DD&E" for rows 2-5, column A (VIN).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "AAAKN3DD&E"
$ws.Range("A3").Value = "AAAKN3DD&E"
$ws.Range("A4").Value = "AAAKN3DD&E"
$ws.Range("A5").Value = "AAAKN3DD&E"

# Move the active selection, matching the recorded cursor position after edit.
$ws.Range("B11").Select() | Out-Null
